# Add a new "Croatia" worksheet at the end of the workbook, based on the
# existing "Spain" sheet (same layout/styles), with market-specific values,
# and make it the active/selected tab (mirroring the "Spain Market" sheet
# that was previously last).

$wb = $excel.ActiveWorkbook

$spain = $wb.Worksheets.Item("Spain")

# Clone the Spain sheet and place the copy right after it.
$spain.Copy([System.Reflection.Missing]::Value, $spain)
$croatia = $wb.Worksheets.Item("Spain (2)")
$croatia.Name = "Croatia"

# Fill in the Croatia-specific values (B4 first, then B2, so the new shared
# strings are appended in that order).
$croatia.Range("B4").Value = "NGC-3103/T2485/T2494"
$croatia.Range("B2").Value = "Croatia Market"

# Spain is no longer the last/active sheet: select the whole sheet (matching
# the other non-active tabs in the workbook) and drop its tab selection.
$spain.Activate()
[void]$spain.Cells.Select()

# Croatia becomes the active tab, with the same kind of cell selection the
# sheet had when it was authored.
$croatia.Activate()
[void]$croatia.Range("D9").Select()
